$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: date 44518 -> 44525; J 50->40; K,L,M 10000->8000; P 667->533
$ws.Range("D2").Value = 44525
$ws.Range("J2").Value = 40
$ws.Range("K2").Value = 8000
$ws.Range("L2").Value = 8000
$ws.Range("M2").Value = 8000
$ws.Range("P2").Value = 533

# Row 3: date 44525 -> 44508; K,L,M 8000->10000; P 533->667 (J stays 40)
$ws.Range("D3").Value = 44508
$ws.Range("K3").Value = 10000
$ws.Range("L3").Value = 10000
$ws.Range("M3").Value = 10000
$ws.Range("P3").Value = 667

# Row 4: date 44508 -> 44518; J 40->50 (K,L,M stay 10000)
$ws.Range("D4").Value = 44518
$ws.Range("J4").Value = 50
